$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-08-28 Wednesday" "2024-08-29 Thursday"

Replace-Text "209÷9=" "563÷4="
Replace-Text "154÷4=" "174÷7="
Replace-Text "616÷6=" "507÷2="
Replace-Text "244÷3=" "733÷9="
Replace-Text "820÷8=" "555÷6="
Replace-Text "194÷8=" "691÷6="
Replace-Text "272÷5=" "852÷3="
Replace-Text "654÷9=" "239÷5="
Replace-Text "374÷7=" "781÷8="
Replace-Text "155÷7=" "640÷9="
Replace-Text "487÷8=" "875÷9="
Replace-Text "370÷6=" "141÷9="
Replace-Text "939÷2=" "318÷8="
Replace-Text "798÷2=" "782÷3="
Replace-Text "666÷7=" "976÷8="
Replace-Text "284÷7=" "850÷5="
Replace-Text "656÷5=" "459÷6="
Replace-Text "104÷3=" "723÷7="
Replace-Text "799÷9=" "708÷3="
Replace-Text "667÷7=" "722÷5="
Replace-Text "574÷6=" "981÷5="
Replace-Text "742÷4=" "950÷2="
Replace-Text "611÷9=" "579÷3="
Replace-Text "975÷8=" "343÷6="
Replace-Text "224÷4=" "170÷9="

Write-Output "done"
